$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.530.41'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.32%  '
$ws.Range("D3").Value = '''1.840.36'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.01%  '
$ws.Range("D4").Value = '''1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").Value = '''259.02'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.49%  '
$ws.Range("E6").Value = '  +0.11%  '
$ws.Range("D7").Value = '''0.5233'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.57%  '
$ws.Range("D8").Value = '''0.3164'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.22%  '
$ws.Range("D9").Value = '''0.06780'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("D10").Value = '''18.71'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.62%  '
$ws.Range("D11").Value = '''0.7794'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.71%  '
$ws.Range("D12").Value = '''0.07781'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.09%  '
$ws.Range("D13").Value = '''1.827.35'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.62%  '
$ws.Range("D14").Value = '''87.91'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.25%  '
$ws.Range("D15").Value = '''5.008'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.42%  '
$ws.Range("D16").Value = '''1.002'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.18%  '
$ws.Range("D17").Value = '''13.84'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.42%  '
$ws.Range("D18").Value = '''1.001'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.06%  '
$ws.Range("D19").Value = '''0.000007929'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.19%  '
$ws.Range("D20").Value = '''26.567.30'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.42%  '
$ws.Range("D21").Value = '''2.067.32'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.15%  '
$ws.Range("D22").Value = '''4.603'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.90%  '
$ws.Range("D23").Value = '''5.963'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("D24").Value = '''9.329'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.51%  '
$ws.Range("D25").Value = '''142.68'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.18%  '
$ws.Range("D26").Value = '''2.215'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.16%  '
$ws.Range("E27").Value = '  +1.34%  '
$ws.Range("E28").Value = '  -0.37%  '
$ws.Range("D29").Value = '''111.99'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.81%  '
$ws.Range("D30").Value = '''4.173'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.00%  '
$ws.Range("E31").Value = '  +0.19%  '
$ws.Range("D32").Value = '''4.066'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.47%  '
$ws.Range("D33").Value = '''0.04887'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.88%  '
$ws.Range("D34").Value = '''1.131'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.53%  '
$ws.Range("D35").Value = '''0.7206'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.86%  '
$ws.Range("D36").Value = '''2.859'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.70%  '
$ws.Range("D37").Value = '''3.095'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.73%  '
$ws.Range("D38").Value = '''2.219'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.49%  '
$ws.Range("D39").Value = '''0.01738'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.80%  '
$ws.Range("D40").Value = '''0.4813'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.23%  '
$ws.Range("E41").Value = '  +1.04%  '
$ws.Range("D42").Value = '''110.23'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.62%  '
$ws.Range("D43").Value = '''5.924'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.39%  '
$ws.Range("E44").Value = '  +0.13%  '
$ws.Range("D45").Value = '''7.638'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.56%  '
$ws.Range("D46").Value = '''0.4159'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.01%  '
$ws.Range("D47").Value = '''8.992'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.57%  '
$ws.Range("D48").Value = '''0.1230'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.03%  '
$ws.Range("D49").Value = '''0.05827'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.71%  '
$ws.Range("D50").Value = '''34.84'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.37%  '
$ws.Range("D51").Value = '''0.8921'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.75%  '
